# Apply "Natmi following Dr Hou advice" update to Il1a-Il1r2 sheet.
# Re-run of the NATMI pipeline added a third sending cluster ("Neutro")
# and refreshed every Ligand/Receptor statistic accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: M1 / Il1a / Il1r2 / ECs
$ws.Cells.Item(2, 1).Value = "M1"
$ws.Cells.Item(2, 2).Value = "Il1a"
$ws.Cells.Item(2, 3).Value = "Il1r2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3.056073
$ws.Cells.Item(2, 8).Value = 9.168219
$ws.Cells.Item(2, 9).Value = 0.4559514113020136
$ws.Cells.Item(2, 10).Value = 0.4559514113020135
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.308173
$ws.Cells.Item(2, 14).Value = 0.6163460000000001
$ws.Cells.Item(2, 15).Value = 0.01562452111805829
$ws.Cells.Item(2, 16).Value = 0.01052496837195314
$ws.Cells.Item(2, 17).Value = 0.9417991846290001
$ws.Cells.Item(2, 18).Value = 5.650795107774001
$ws.Cells.Item(2, 19).Value = 0.007124022454696791
$ws.Cells.Item(2, 20).Value = 0.00479887418310109

# Row 3: M1 / Il1a / Il1r2 / FAPs
$ws.Cells.Item(3, 1).Value = "M1"
$ws.Cells.Item(3, 2).Value = "Il1a"
$ws.Cells.Item(3, 3).Value = "Il1r2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3.056073
$ws.Cells.Item(3, 8).Value = 9.168219
$ws.Cells.Item(3, 9).Value = 0.4559514113020136
$ws.Cells.Item(3, 10).Value = 0.4559514113020135
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.770252
$ws.Cells.Item(3, 14).Value = 5.310756
$ws.Cells.Item(3, 15).Value = 0.08975263815546759
$ws.Cells.Item(3, 16).Value = 0.09068857254068391
$ws.Cells.Item(3, 17).Value = 5.410019340396
$ws.Cells.Item(3, 18).Value = 48.690174063564
$ws.Cells.Item(3, 19).Value = 0.0409228420350644
$ws.Cells.Item(3, 20).Value = 0.04134958263888985

# Row 4: M1 / Il1a / Il1r2 / M1
$ws.Cells.Item(4, 1).Value = "M1"
$ws.Cells.Item(4, 2).Value = "Il1a"
$ws.Cells.Item(4, 3).Value = "Il1r2"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 3.056073
$ws.Cells.Item(4, 8).Value = 9.168219
$ws.Cells.Item(4, 9).Value = 0.4559514113020136
$ws.Cells.Item(4, 10).Value = 0.4559514113020135
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 10.33148
$ws.Cells.Item(4, 14).Value = 30.99444
$ws.Cells.Item(4, 15).Value = 0.5238110653457533
$ws.Cells.Item(4, 16).Value = 0.5292733313859411
$ws.Cells.Item(4, 17).Value = 31.57375707804
$ws.Cells.Item(4, 18).Value = 284.16381370236
$ws.Cells.Item(4, 19).Value = 0.2388323945000075
$ws.Cells.Item(4, 20).Value = 0.2413229224099381

# Row 5: M1 / Il1a / Il1r2 / M2
$ws.Cells.Item(5, 1).Value = "M1"
$ws.Cells.Item(5, 2).Value = "Il1a"
$ws.Cells.Item(5, 3).Value = "Il1r2"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.056073
$ws.Cells.Item(5, 8).Value = 9.168219
$ws.Cells.Item(5, 9).Value = 0.4559514113020136
$ws.Cells.Item(5, 10).Value = 0.4559514113020135
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 6.856480333333334
$ws.Cells.Item(5, 14).Value = 20.569441
$ws.Cells.Item(5, 15).Value = 0.3476268906222089
$ws.Cells.Item(5, 16).Value = 0.3512519201126578
$ws.Cells.Item(5, 17).Value = 20.953904421731
$ws.Cells.Item(5, 18).Value = 188.585139795579
$ws.Cells.Item(5, 19).Value = 0.1585009713857269
$ws.Cells.Item(5, 20).Value = 0.1601538086979084

# Row 6: M1 / Il1a / Il1r2 / Neutro
$ws.Cells.Item(6, 1).Value = "M1"
$ws.Cells.Item(6, 2).Value = "Il1a"
$ws.Cells.Item(6, 3).Value = "Il1r2"
$ws.Cells.Item(6, 4).Value = "Neutro"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 3.056073
$ws.Cells.Item(6, 8).Value = 9.168219
$ws.Cells.Item(6, 9).Value = 0.4559514113020136
$ws.Cells.Item(6, 10).Value = 0.4559514113020135
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.1548006666666667
$ws.Cells.Item(6, 14).Value = 0.464402
$ws.Cells.Item(6, 15).Value = 0.007848469156684183
$ws.Cells.Item(6, 16).Value = 0.007930312457405064
$ws.Cells.Item(6, 17).Value = 0.473082137782
$ws.Cells.Item(6, 18).Value = 4.257739240038
$ws.Cells.Item(6, 19).Value = 0.003578520588550477
$ws.Cells.Item(6, 20).Value = 0.003615837157019778

# Row 7: M1 / Il1a / Il1r2 / sCs
$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Il1a"
$ws.Cells.Item(7, 3).Value = "Il1r2"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 3.056073
$ws.Cells.Item(7, 8).Value = 9.168219
$ws.Cells.Item(7, 9).Value = 0.4559514113020136
$ws.Cells.Item(7, 10).Value = 0.4559514113020135
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3024905
$ws.Cells.Item(7, 14).Value = 0.604981
$ws.Cells.Item(7, 15).Value = 0.01533641560182758
$ws.Cells.Item(7, 16).Value = 0.01033089513135898
$ws.Cells.Item(7, 17).Value = 0.9244330498065
$ws.Cells.Item(7, 18).Value = 5.546598298839
$ws.Cells.Item(7, 19).Value = 0.006992660337967504
$ws.Cells.Item(7, 20).Value = 0.004710386215156227

# Row 8: M2 / Il1a / Il1r2 / ECs
$ws.Cells.Item(8, 1).Value = "M2"
$ws.Cells.Item(8, 2).Value = "Il1a"
$ws.Cells.Item(8, 3).Value = "Il1r2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 3.592139
$ws.Cells.Item(8, 8).Value = 10.776417
$ws.Cells.Item(8, 9).Value = 0.5359298834298145
$ws.Cells.Item(8, 10).Value = 0.5359298834298145
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.308173
$ws.Cells.Item(8, 14).Value = 0.6163460000000001
$ws.Cells.Item(8, 15).Value = 0.01562452111805829
$ws.Cells.Item(8, 16).Value = 0.01052496837195314
$ws.Cells.Item(8, 17).Value = 1.107000252047
$ws.Cells.Item(8, 18).Value = 6.642001512282001
$ws.Cells.Item(8, 19).Value = 0.008373647781447653
$ws.Cells.Item(8, 20).Value = 0.005640645072683331

# Row 9: M2 / Il1a / Il1r2 / FAPs
$ws.Cells.Item(9, 1).Value = "M2"
$ws.Cells.Item(9, 2).Value = "Il1a"
$ws.Cells.Item(9, 3).Value = "Il1r2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 3.592139
$ws.Cells.Item(9, 8).Value = 10.776417
$ws.Cells.Item(9, 9).Value = 0.5359298834298145
$ws.Cells.Item(9, 10).Value = 0.5359298834298145
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.770252
$ws.Cells.Item(9, 14).Value = 5.310756
$ws.Cells.Item(9, 15).Value = 0.08975263815546759
$ws.Cells.Item(9, 16).Value = 0.09068857254068391
$ws.Cells.Item(9, 17).Value = 6.358991249028
$ws.Cells.Item(9, 18).Value = 57.230921241252
$ws.Cells.Item(9, 19).Value = 0.04810112090417807
$ws.Cells.Item(9, 20).Value = 0.04860271611014501

# Row 10: M2 / Il1a / Il1r2 / M1
$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Il1a"
$ws.Cells.Item(10, 3).Value = "Il1r2"
$ws.Cells.Item(10, 4).Value = "M1"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 3.592139
$ws.Cells.Item(10, 8).Value = 10.776417
$ws.Cells.Item(10, 9).Value = 0.5359298834298145
$ws.Cells.Item(10, 10).Value = 0.5359298834298145
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 10.33148
$ws.Cells.Item(10, 14).Value = 30.99444
$ws.Cells.Item(10, 15).Value = 0.5238110653457533
$ws.Cells.Item(10, 16).Value = 0.5292733313859411
$ws.Cells.Item(10, 17).Value = 37.11211223572001
$ws.Cells.Item(10, 18).Value = 334.00901012148
$ws.Cells.Item(10, 19).Value = 0.2807260031899965
$ws.Cells.Item(10, 20).Value = 0.283653394792177

# Row 11: M2 / Il1a / Il1r2 / M2
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Il1a"
$ws.Cells.Item(11, 3).Value = "Il1r2"
$ws.Cells.Item(11, 4).Value = "M2"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 3.592139
$ws.Cells.Item(11, 8).Value = 10.776417
$ws.Cells.Item(11, 9).Value = 0.5359298834298145
$ws.Cells.Item(11, 10).Value = 0.5359298834298145
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 6.856480333333334
$ws.Cells.Item(11, 14).Value = 20.569441
$ws.Cells.Item(11, 15).Value = 0.3476268906222089
$ws.Cells.Item(11, 16).Value = 0.3512519201126578
$ws.Cells.Item(11, 17).Value = 24.62943040809967
$ws.Cells.Item(11, 18).Value = 221.664873672897
$ws.Cells.Item(11, 19).Value = 0.1863036389682293
$ws.Cells.Item(11, 20).Value = 0.1882464006004752

# Row 12: M2 / Il1a / Il1r2 / Neutro
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Il1a"
$ws.Cells.Item(12, 3).Value = "Il1r2"
$ws.Cells.Item(12, 4).Value = "Neutro"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 3.592139
$ws.Cells.Item(12, 8).Value = 10.776417
$ws.Cells.Item(12, 9).Value = 0.5359298834298145
$ws.Cells.Item(12, 10).Value = 0.5359298834298145
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.1548006666666667
$ws.Cells.Item(12, 14).Value = 0.464402
$ws.Cells.Item(12, 15).Value = 0.007848469156684183
$ws.Cells.Item(12, 16).Value = 0.007930312457405064
$ws.Cells.Item(12, 17).Value = 0.5560655119593333
$ws.Cells.Item(12, 18).Value = 5.004589607634
$ws.Cells.Item(12, 19).Value = 0.004206229160244249
$ws.Cells.Item(12, 20).Value = 0.004250091430859102

# Row 13: M2 / Il1a / Il1r2 / sCs
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Il1a"
$ws.Cells.Item(13, 3).Value = "Il1r2"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 3.592139
$ws.Cells.Item(13, 8).Value = 10.776417
$ws.Cells.Item(13, 9).Value = 0.5359298834298145
$ws.Cells.Item(13, 10).Value = 0.5359298834298145
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.3024905
$ws.Cells.Item(13, 14).Value = 0.604981
$ws.Cells.Item(13, 15).Value = 0.01533641560182758
$ws.Cells.Item(13, 16).Value = 0.01033089513135898
$ws.Cells.Item(13, 17).Value = 1.0865879221795
$ws.Cells.Item(13, 18).Value = 6.519527533077
$ws.Cells.Item(13, 19).Value = 0.008219243425718642
$ws.Cells.Item(13, 20).Value = 0.005536635423474857

# Row 14: Neutro / Il1a / Il1r2 / ECs
$ws.Cells.Item(14, 1).Value = "Neutro"
$ws.Cells.Item(14, 2).Value = "Il1a"
$ws.Cells.Item(14, 3).Value = "Il1r2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.05441666666666667
$ws.Cells.Item(14, 8).Value = 0.16325
$ws.Cells.Item(14, 9).Value = 0.0081187052681719
$ws.Cells.Item(14, 10).Value = 0.008118705268171898
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.308173
$ws.Cells.Item(14, 14).Value = 0.6163460000000001
$ws.Cells.Item(14, 15).Value = 0.01562452111805829
$ws.Cells.Item(14, 16).Value = 0.01052496837195314
$ws.Cells.Item(14, 17).Value = 0.01676974741666667
$ws.Cells.Item(14, 18).Value = 0.1006184845
$ws.Cells.Item(14, 19).Value = 0.0001268508819138429
$ws.Cells.Item(14, 20).Value = 0.00008544911616871857

# Row 15: Neutro / Il1a / Il1r2 / FAPs
$ws.Cells.Item(15, 1).Value = "Neutro"
$ws.Cells.Item(15, 2).Value = "Il1a"
$ws.Cells.Item(15, 3).Value = "Il1r2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.05441666666666667
$ws.Cells.Item(15, 8).Value = 0.16325
$ws.Cells.Item(15, 9).Value = 0.0081187052681719
$ws.Cells.Item(15, 10).Value = 0.008118705268171898
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.770252
$ws.Cells.Item(15, 14).Value = 5.310756
$ws.Cells.Item(15, 15).Value = 0.08975263815546759
$ws.Cells.Item(15, 16).Value = 0.09068857254068391
$ws.Cells.Item(15, 17).Value = 0.096331213
$ws.Cells.Item(15, 18).Value = 0.866980917
$ws.Cells.Item(15, 19).Value = 0.0007286752162251211
$ws.Cells.Item(15, 20).Value = 0.0007362737916490398

# Row 16: Neutro / Il1a / Il1r2 / M1
$ws.Cells.Item(16, 1).Value = "Neutro"
$ws.Cells.Item(16, 2).Value = "Il1a"
$ws.Cells.Item(16, 3).Value = "Il1r2"
$ws.Cells.Item(16, 4).Value = "M1"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.05441666666666667
$ws.Cells.Item(16, 8).Value = 0.16325
$ws.Cells.Item(16, 9).Value = 0.0081187052681719
$ws.Cells.Item(16, 10).Value = 0.008118705268171898
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 10.33148
$ws.Cells.Item(16, 14).Value = 30.99444
$ws.Cells.Item(16, 15).Value = 0.5238110653457533
$ws.Cells.Item(16, 16).Value = 0.5292733313859411
$ws.Cells.Item(16, 17).Value = 0.5622047033333334
$ws.Cells.Item(16, 18).Value = 5.05984233
$ws.Cells.Item(16, 19).Value = 0.004252667655749303
$ws.Cells.Item(16, 20).Value = 0.004297014183825931

# Row 17: Neutro / Il1a / Il1r2 / M2
$ws.Cells.Item(17, 1).Value = "Neutro"
$ws.Cells.Item(17, 2).Value = "Il1a"
$ws.Cells.Item(17, 3).Value = "Il1r2"
$ws.Cells.Item(17, 4).Value = "M2"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.05441666666666667
$ws.Cells.Item(17, 8).Value = 0.16325
$ws.Cells.Item(17, 9).Value = 0.0081187052681719
$ws.Cells.Item(17, 10).Value = 0.008118705268171898
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 6.856480333333334
$ws.Cells.Item(17, 14).Value = 20.569441
$ws.Cells.Item(17, 15).Value = 0.3476268906222089
$ws.Cells.Item(17, 16).Value = 0.3512519201126578
$ws.Cells.Item(17, 17).Value = 0.3731068048055556
$ws.Cells.Item(17, 18).Value = 3.35796124325
$ws.Cells.Item(17, 19).Value = 0.002822280268252745
$ws.Cells.Item(17, 20).Value = 0.002851710814274129

# Row 18: Neutro / Il1a / Il1r2 / Neutro
$ws.Cells.Item(18, 1).Value = "Neutro"
$ws.Cells.Item(18, 2).Value = "Il1a"
$ws.Cells.Item(18, 3).Value = "Il1r2"
$ws.Cells.Item(18, 4).Value = "Neutro"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.05441666666666667
$ws.Cells.Item(18, 8).Value = 0.16325
$ws.Cells.Item(18, 9).Value = 0.0081187052681719
$ws.Cells.Item(18, 10).Value = 0.008118705268171898
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.1548006666666667
$ws.Cells.Item(18, 14).Value = 0.464402
$ws.Cells.Item(18, 15).Value = 0.007848469156684183
$ws.Cells.Item(18, 16).Value = 0.007930312457405064
$ws.Cells.Item(18, 17).Value = 0.008423736277777778
$ws.Cells.Item(18, 18).Value = 0.0758136265
$ws.Cells.Item(18, 19).Value = 0.00006371940788945655
$ws.Cells.Item(18, 20).Value = 0.00006438386952618372

# Row 19: Neutro / Il1a / Il1r2 / sCs
$ws.Cells.Item(19, 1).Value = "Neutro"
$ws.Cells.Item(19, 2).Value = "Il1a"
$ws.Cells.Item(19, 3).Value = "Il1r2"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.05441666666666667
$ws.Cells.Item(19, 8).Value = 0.16325
$ws.Cells.Item(19, 9).Value = 0.0081187052681719
$ws.Cells.Item(19, 10).Value = 0.008118705268171898
$ws.Cells.Item(19, 11).Value = 2
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.3024905
$ws.Cells.Item(19, 14).Value = 0.604981
$ws.Cells.Item(19, 15).Value = 0.01533641560182758
$ws.Cells.Item(19, 16).Value = 0.01033089513135898
$ws.Cells.Item(19, 17).Value = 0.01646052470833333
$ws.Cells.Item(19, 18).Value = 0.09876314825
$ws.Cells.Item(19, 19).Value = 0.0001245118381414313
$ws.Cells.Item(19, 20).Value = 0.00008387349272789557

